$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '98.531.68'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.342.93'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.70'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '647.30'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +10.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.458'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +15.68%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +22.25%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.339.55'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '44.08'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +21.47%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000269'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +7.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '98.071.33'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.973.38'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.55'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.335.31'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.41'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +19.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.77'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +9.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '533.45'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +7.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.56'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.18'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +8.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000213'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.28%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +50.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '103.20'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +14.05%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +7.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.66'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.523.55'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.151'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +12.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +13.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.188'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '29.03'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.520'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +10.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.73'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.156'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.06'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '513.19'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.62%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.76'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.32'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.87%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.39%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0408'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +25.20%  '
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '163.90'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.68'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +16.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.94'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.06%  '
